$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "59.297.60"
$ws.Range("E2").Value = "  -5.38%  "
Set-TextValue "D3" "2.483.32"
$ws.Range("E3").Value = "  -7.07%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.24%  "
Set-TextValue "D5" "541.30"
$ws.Range("E5").Value = "  -2.11%  "
Set-TextValue "D6" "146.28"
$ws.Range("E6").Value = "  -6.62%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  -1.99%  "
Set-TextValue "D9" "2.480.26"
$ws.Range("E9").Value = "  -7.24%  "
Set-TextValue "D10" "0.100"
$ws.Range("E10").Value = "  -4.40%  "
$ws.Range("E11").Value = "  -1.05%  "
Set-TextValue "D12" "5.52"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("E13").Value = "  -2.96%  "
Set-TextValue "D14" "2.900.00"
$ws.Range("E14").Value = "  -7.87%  "
Set-TextValue "D15" "24.38"
$ws.Range("E15").Value = "  -6.21%  "
Set-TextValue "D16" "59.257.02"
$ws.Range("E16").Value = "  -5.37%  "
$ws.Range("E17").Value = "  -4.23%  "
Set-TextValue "D18" "2.479.02"
$ws.Range("E18").Value = "  -7.47%  "
$ws.Range("E19").Value = "  -4.59%  "
Set-TextValue "D20" "4.39"
$ws.Range("E20").Value = "  -3.81%  "
Set-TextValue "D21" "324.18"
$ws.Range("E21").Value = "  -5.45%  "
$ws.Range("E22").Value = "  -1.76%  "
Set-TextValue "D23" "5.79"
$ws.Range("E23").Value = "  -6.11%  "
Set-TextValue "D24" "60.93"
$ws.Range("E24").Value = "  -3.43%  "
Set-TextValue "D25" "0.454"
$ws.Range("E25").Value = "  -10.16%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D26" "0.993"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D27" "0.161"
$ws.Range("E27").Value = "  -4.76%  "
Set-TextValue "D28" "7.75"
$ws.Range("E28").Value = "  -4.29%  "
$ws.Range("E29").Value = "  -5.70%  "
$ws.Range("E30").Value = "  -8.33%  "
Set-TextValue "D31" "0.0₃0779"
$ws.Range("E31").Value = "  -7.65%  "
Set-TextValue "D32" "1.82"
$ws.Range("E32").Value = "  -4.98%  "
Set-TextValue "D34" "158.10"
$ws.Range("E34").Value = "  -2.48%  "
Set-TextValue "D35" "1.42"
$ws.Range("E35").Value = "  -1.51%  "
Set-TextValue "D36" "18.70"
$ws.Range("E36").Value = "  -3.37%  "
$ws.Range("E37").Value = "  -7.00%  "
Set-TextValue "D38" "1.73"
$ws.Range("E38").Value = "  -2.29%  "
Set-TextValue "D39" "5.89"
$ws.Range("E39").Value = "  -3.89%  "
Set-TextValue "D40" "313.94"
$ws.Range("E40").Value = "  -6.40%  "
Set-TextValue "D41" "36.55"
$ws.Range("E41").Value = "  -4.45%  "
$ws.Range("E42").Value = "  -5.61%  "
Set-TextValue "D43" "0.833"
$ws.Range("E43").Value = "  -9.12%  "
Set-TextValue "D44" "0.997"
$ws.Range("E44").Value = "  -0.22%  "
Set-TextValue "D45" "0.601"
$ws.Range("E45").Value = "  -2.07%  "
Set-TextValue "D46" "10.72"
$ws.Range("E46").Value = "  -2.63%  "
Set-TextValue "D47" "126.80"
$ws.Range("E47").Value = "  -1.96%  "
Set-TextValue "D48" "0.0532"
$ws.Range("E48").Value = "  -3.75%  "
Set-TextValue "D49" "0.0937"
$ws.Range("E49").Value = "  -2.97%  "
$ws.Range("E50").Value = "  -3.13%  "
Set-TextValue "D51" "18.54"
$ws.Range("E51").Value = "  -7.50%  "
